# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The worker/period detail table (Hoja1!B16:G35) is re-sorted: instead of
# being grouped by period-then-worker, it is regrouped by worker-then-period
# (periods newest-first per worker). The underlying facts (doc type, doc
# number, worker name, period, "Valor Mora" and "Salario Basico") are the
# same set of records - only their row order changes - except the Feb/2025
# ("2501") "Valor Mora" figure, which is corrected from 52000 to 50266 for
# every worker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: DocType, DocNumber, WorkerName, Period, ValorMora, SalarioBasico
$rows = @(
    @("CC", "9158235",  "FERNANDO JOSE JIMENEZ BARRETO",     "2501", 50266, 1300000),
    @("CC", "9158235",  "FERNANDO JOSE JIMENEZ BARRETO",     "2412", 52000, 1300000),
    @("CC", "9158235",  "FERNANDO JOSE JIMENEZ BARRETO",     "2409", 52000, 1300000),
    @("CC", "9158235",  "FERNANDO JOSE JIMENEZ BARRETO",     "2408", 52000, 1300000),
    @("CC", "32907838", "BANIDIS DEL CARMEN BANQUET BLANCO", "2501", 50266, 1300000),
    @("CC", "32907838", "BANIDIS DEL CARMEN BANQUET BLANCO", "2412", 52000, 1300000),
    @("CC", "32907838", "BANIDIS DEL CARMEN BANQUET BLANCO", "2411", 52000, 1300000),
    @("CC", "32907838", "BANIDIS DEL CARMEN BANQUET BLANCO", "2410", 52000, 1300000),
    @("CC", "32907838", "BANIDIS DEL CARMEN BANQUET BLANCO", "2409", 52000, 1300000),
    @("CC", "32907838", "BANIDIS DEL CARMEN BANQUET BLANCO", "2408", 52000, 1300000),
    @("CC", "45579779", "GLADYS DEL SOCORRO LEGUIA ROBLES",  "2501", 50266, 1300000),
    @("CC", "45579779", "GLADYS DEL SOCORRO LEGUIA ROBLES",  "2412", 52000, 1300000),
    @("CC", "45579779", "GLADYS DEL SOCORRO LEGUIA ROBLES",  "2411", 52000, 1300000),
    @("CC", "45579779", "GLADYS DEL SOCORRO LEGUIA ROBLES",  "2410", 52000, 1300000),
    @("CC", "45579779", "GLADYS DEL SOCORRO LEGUIA ROBLES",  "2409", 52000, 1300000),
    @("CC", "45579779", "GLADYS DEL SOCORRO LEGUIA ROBLES",  "2408", 52000, 1300000),
    @("CC", "73121215", "YIDIO MORALES TORRES",              "2501", 50266, 1300000),
    @("CC", "73121215", "YIDIO MORALES TORRES",              "2412", 52000, 1300000),
    @("CC", "73121215", "YIDIO MORALES TORRES",              "2409", 52000, 1300000),
    @("CC", "73121215", "YIDIO MORALES TORRES",              "2408", 52000, 1300000)
)

$firstRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $firstRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 2).Value = $data[0]   # B - Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value = $data[1]   # C - N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $data[2]   # D - Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $data[3]   # E - Periodo Mora
    $ws.Cells.Item($r, 6).Value = $data[4]   # F - Valor Mora
    $ws.Cells.Item($r, 7).Value = $data[5]   # G - Salario Basico
}
